# Apply "Models in file updated" changes to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated existing values in the Models table (columns A-I) ---
$ws.Range("A2").Value = 3
$ws.Range("E2").Value = 50
$ws.Range("K2").Value = 1

$ws.Range("A3").Value = 6
$ws.Range("C3").Value = 97
$ws.Range("E3").Value = 50
$ws.Range("M3").Value = 630

$ws.Range("A4").Value = 8
$ws.Range("C4").Value = 99
$ws.Range("E4").Value = 51

$ws.Range("A6").Value = 8
$ws.Range("C6").Value = 91
$ws.Range("E6").Value = 51

$ws.Range("A7").Value = 6
$ws.Range("E7").Value = 50

# --- New L/M columns with QTY / MODEL breakdown for rows 11-24 ---
$lmData = @(
    @(11, 3, "VTY"),
    @(12, 1, "VT"),
    @(13, 4, "VTL"),
    @(14, 2, "VTW"),
    @(15, 1, "VTY"),
    @(16, 1, "VTY"),
    @(17, 2, "VTW"),
    @(18, 2, "VTW"),
    @(19, 1, "VT"),
    @(20, 1, "VTY"),
    @(21, 1, "VTW"),
    @(22, 1, "VTW"),
    @(23, 1, "VTY"),
    @(24, 1, "VTY")
)

foreach ($entry in $lmData) {
    $row = $entry[0]
    $qty = $entry[1]
    $model = $entry[2]
    $ws.Cells.Item($row, 12).Value = $qty
    $ws.Cells.Item($row, 13).Value = $model
}

# Update the active selection to match the saved state of the sheet.
$ws.Range("A7").Select()
